$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "              -Decision Table Testing" -> "             -Decision Table Testing"
#    (one fewer leading space before the dash that precedes "Decision Table Testing")
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("              -Decision Table Testing", $true, $true, $false, $false, $false, $true, 1, $false, "             -Decision Table Testing", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) The empty paragraph right after "Black Box Testing" / "-Decision Table Testing"
#    becomes a numbered ("Black Box Testing" list, numId 1) list item containing "OOPS",
#    and a brand-new "List Paragraph" styled paragraph containing "             -Inheritance"
#    is inserted immediately after it (right before the following table).
# ---------------------------------------------------------------------------
$target = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "`r") {
        if ($cand.Style.NameLocal -eq "Normal") {
            $prev = $d.Paragraphs.Item($i - 1).Range.Text
            if ($prev -like "*Decision Table Testing*") {
                $target = $cand
                break
            }
        }
    }
}

$r = $target.Range
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr><w:t>OOPS</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="404040" w:themeColor="text1" w:themeTint="BF"/></w:rPr><w:t xml:space="preserve">             -Inheritance</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
